# Auto-generated edit script: apply numeric updates to the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit's scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 450.22223
$ws.Range("I2").Value = 419
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 419
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -306
$ws.Range("N2").Value = -926

$ws.Range("H137").Value = 531452.0600000001
$ws.Range("I137").Value = 2609.5386
$ws.Range("K137").Value = 7828.6158
$ws.Range("M137").Value = -5278.6158

$ws.Range("H138").Value = 4663.1274
$ws.Range("I138").Value = 1834.4736
$ws.Range("J138").Value = 6582.5713
$ws.Range("K138").Value = 5503.4208
$ws.Range("L138").Value = 19747.7139
$ws.Range("M138").Value = -363.4207999999999
$ws.Range("N138").Value = -30027.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12302.4375
$ws.Range("I32").Value = 12710.11
$ws.Range("K32").Value = 12710.11
$ws.Range("M32").Value = -12423.11

$ws.Range("H45").Value = 1756.0385
$ws.Range("I45").Value = 1779.4584
$ws.Range("J45").Value = 1475
$ws.Range("K45").Value = 1779.4584
$ws.Range("L45").Value = 1475
$ws.Range("M45").Value = -1402.4584
$ws.Range("N45").Value = -2229

$ws.Range("H61").Value = 5803.132
$ws.Range("I61").Value = 2205.625
$ws.Range("J61").Value = 16872.385
$ws.Range("K61").Value = 2205.625
$ws.Range("L61").Value = 16872.385
$ws.Range("M61").Value = -1993.625
$ws.Range("N61").Value = -17296.385

$ws.Range("H74").Value = 4500.421
$ws.Range("I74").Value = 1701.8462
$ws.Range("J74").Value = 10564
$ws.Range("K74").Value = 1701.8462
$ws.Range("L74").Value = 10564
$ws.Range("M74").Value = -827.8462
$ws.Range("N74").Value = -12312

$ws.Range("H77").Value = 4500.421
$ws.Range("I77").Value = 1701.8462
$ws.Range("J77").Value = 10564
$ws.Range("K77").Value = 8509.231
$ws.Range("L77").Value = 52820
$ws.Range("M77").Value = -4141.231
$ws.Range("N77").Value = -61556

$ws.Range("H97").Value = 1504.3572
$ws.Range("I97").Value = 1512.5
$ws.Range("J97").Value = 1493.5
$ws.Range("K97").Value = 1512.5
$ws.Range("L97").Value = 1493.5
$ws.Range("M97").Value = -1016.5
$ws.Range("N97").Value = -2485.5

$ws.Range("H136").Value = 5803.132
$ws.Range("I136").Value = 2205.625
$ws.Range("J136").Value = 16872.385
$ws.Range("K136").Value = 6616.875
$ws.Range("L136").Value = 50617.155
$ws.Range("M136").Value = -4066.875
$ws.Range("N136").Value = -55717.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27018.77
$ws.Range("I134").Value = 1490.9697
$ws.Range("J134").Value = 167421.67
$ws.Range("K134").Value = 4472.909100000001
$ws.Range("L134").Value = 502265.01
$ws.Range("M134").Value = -1937.909100000001
$ws.Range("N134").Value = -507335.01

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 225.5
$ws.Range("I22").Value = 225.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 225.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 124.5
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 546050.5
$ws.Range("I31").Value = 12448.059
$ws.Range("J31").Value = 731177.9399999999
$ws.Range("K31").Value = 12448.059
$ws.Range("L31").Value = 731177.9399999999
$ws.Range("M31").Value = -12153.059
$ws.Range("N31").Value = -731767.9399999999

$ws.Range("H34").Value = 546050.5
$ws.Range("I34").Value = 12448.059
$ws.Range("J34").Value = 731177.9399999999
$ws.Range("K34").Value = 12448.059
$ws.Range("L34").Value = 731177.9399999999
$ws.Range("M34").Value = -12246.059
$ws.Range("N34").Value = -731581.9399999999

$ws.Range("H122").Value = 3903.2144
$ws.Range("I122").Value = 3447.889
$ws.Range("J122").Value = 4722.8
$ws.Range("K122").Value = 10343.667
$ws.Range("L122").Value = 14168.4
$ws.Range("M122").Value = -7893.667000000001
$ws.Range("N122").Value = -19068.4

$ws.Range("H132").Value = 3578.5557
$ws.Range("I132").Value = 3314.923
$ws.Range("J132").Value = 4264
$ws.Range("K132").Value = 9944.769
$ws.Range("L132").Value = 12792
$ws.Range("M132").Value = -7414.769
$ws.Range("N132").Value = -17852

$ws.Range("H134").Value = 1815.5
$ws.Range("I134").Value = 1393.1459
$ws.Range("J134").Value = 3842.8
$ws.Range("K134").Value = 4179.4377
$ws.Range("L134").Value = 11528.4
$ws.Range("M134").Value = -1644.4377
$ws.Range("N134").Value = -16598.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5559466
$ws.Range("I5").Value = 435
$ws.Range("J5").Value = 6949223.5
$ws.Range("K5").Value = 1305
$ws.Range("L5").Value = 20847670.5
$ws.Range("M5").Value = -1193
$ws.Range("N5").Value = -20847894.5

$ws.Range("H132").Value = 4249.75
$ws.Range("J132").Value = 3428.2856
$ws.Range("L132").Value = 30854.5704
$ws.Range("N132").Value = -35914.5704

$ws.Range("H135").Value = 5559466
$ws.Range("I135").Value = 435
$ws.Range("J135").Value = 6949223.5
$ws.Range("K135").Value = 3915
$ws.Range("L135").Value = 62543011.5
$ws.Range("M135").Value = -1380
$ws.Range("N135").Value = -62548081.5

$ws.Range("H137").Value = 30678
$ws.Range("I137").Value = 1722.5333
$ws.Range("J137").Value = 175455.33
$ws.Range("K137").Value = 5167.5999
$ws.Range("L137").Value = 526365.99
$ws.Range("M137").Value = -67.59990000000016
$ws.Range("N137").Value = -536565.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3240.423
$ws.Range("I102").Value = 2678.0588
$ws.Range("J102").Value = 4302.6665
$ws.Range("K102").Value = 2678.0588
$ws.Range("L102").Value = 4302.6665
$ws.Range("M102").Value = -1056.0588
$ws.Range("N102").Value = -7546.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6252.5835
$ws.Range("I122").Value = 6089.6
$ws.Range("J122").Value = 7067.5
$ws.Range("K122").Value = 18268.8
$ws.Range("L122").Value = 21202.5
$ws.Range("M122").Value = -15818.8
$ws.Range("N122").Value = -26102.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1013.2857
$ws.Range("I122").Value = 1003.2727
$ws.Range("J122").Value = 1050
$ws.Range("K122").Value = 3009.8181
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -559.8181
$ws.Range("N122").Value = -8050

$ws.Range("H123").Value = 50417.668
$ws.Range("J123").Value = 50417.668
$ws.Range("L123").Value = 50417.668
$ws.Range("N123").Value = -60217.668

$ws.Range("H126").Value = 1615.56
$ws.Range("I126").Value = 1627.2
$ws.Range("J126").Value = 1569
$ws.Range("K126").Value = 4881.6
$ws.Range("L126").Value = 4707
$ws.Range("M126").Value = -2411.6
$ws.Range("N126").Value = -9647

Write-Host "Applied scheduled-runner refresh to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve profit sheets."
